# Updated cryptos list values (price + volume(1h) columns), plus two rank swaps
# (ARBITRUM/MXToken and TrustWalletToken/Quant exchanging row positions).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name (B) updates: two pairs of rows swapped rank order ---
$bUpdates = @{
    "B36" = "ARBITRUM"
    "B37" = "MXToken"
    "B43" = "TrustWalletToken"
    "B44" = "Quant"
}

# --- Link (C) updates: follow the B swaps ---
$cUpdates = @{
    "C36" = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
    "C37" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "C43" = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
    "C44" = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
}

# --- Price (D) updates. These are plain display text (e.g. "25.869.21",
# "1.001", "0.0₈117") and must stay text, never get reinterpreted as
# numbers/dates, so NumberFormat is forced to "@" (Text) before assignment. ---
$dUpdates = @{
    "D2" = "25.869.21"
    "D3" = "1.637.38"
    "D4" = "1.001"
    "D5" = "215.27"
    "D6" = "0.5037"
    "D7" = "1.003"
    "D8" = "0.2565"
    "D9" = "0.06382"
    "D10" = "19.67"
    "D11" = "0.07728"
    "D12" = "1.651.49"
    "D13" = "4.251"
    "D14" = "1.864.04"
    "D15" = "0.5459"
    "D17" = "64.07"
    "D18" = "25.894.82"
    "D19" = "1.002"
    "D22" = "9.883"
    "D23" = "5.963"
    "D24" = "1.004"
    "D25" = "1.911"
    "D26" = "140.78"
    "D28" = "15.64"
    "D29" = "6.750"
    "D30" = "1.244"
    "D31" = "0.04964"
    "D32" = "3.272"
    "D35" = "2.376"
    "D36" = "0.8934"
    "D37" = "2.623"
    "D38" = "1.152.45"
    "D39" = "0.5605"
    "D41" = "1.004"
    "D42" = "5.655"
    "D43" = "0.8072"
    "D44" = "99.80"
    "D45" = "1.775.58"
    "D46" = "0.0₈117"
    "D47" = "0.4527"
    "D49" = "54.79"
    "D50" = "0.05052"
}

# --- Volume(1h) (E) updates: always padded/percent strings, safe as plain text ---
$eUpdates = @{
    "E2" = "  -1.71%  "
    "E3" = "  -1.44%  "
    "E4" = "  -0.70%  "
    "E5" = "  -0.97%  "
    "E6" = "  -2.13%  "
    "E7" = "  -0.70%  "
    "E8" = "  -1.60%  "
    "E9" = "  -1.51%  "
    "E10" = "  -1.69%  "
    "E11" = "  -1.46%  "
    "E12" = "  -0.66%  "
    "E13" = "  -1.54%  "
    "E14" = "  -1.33%  "
    "E15" = "  -1.85%  "
    "E16" = "  -2.27%  "
    "E17" = "  -0.76%  "
    "E18" = "  -1.54%  "
    "E19" = "  -0.55%  "
    "E20" = "  -3.93%  "
    "E21" = "  -1.36%  "
    "E22" = "  -2.44%  "
    "E23" = "  -1.43%  "
    "E24" = "  -0.57%  "
    "E25" = "  +8.36%  "
    "E26" = "  -2.87%  "
    "E27" = "  -3.64%  "
    "E28" = "  -1.74%  "
    "E29" = "  -3.89%  "
    "E30" = "  -0.76%  "
    "E31" = "  -3.15%  "
    "E32" = "  -3.11%  "
    "E33" = "  -1.62%  "
    "E34" = "  -1.68%  "
    "E35" = "  +0.54%  "
    "E37" = "  -4.50%  "
    "E38" = "  -1.76%  "
    "E39" = "  -2.62%  "
    "E40" = "  -1.78%  "
    "E41" = "  -0.73%  "
    "E42" = "  -1.41%  "
    "E43" = "  -2.55%  "
    "E44" = "  -0.97%  "
    "E45" = "  -1.27%  "
    "E46" = "  +6.36%  "
    "E47" = "  -0.78%  "
    "E48" = "  -0.64%  "
    "E49" = "  -1.64%  "
    "E50" = "  -0.82%  "
    "E51" = "  -0.49%  "
}

foreach ($ref in $bUpdates.Keys) {
    $ws.Range($ref).Value = $bUpdates[$ref]
}

foreach ($ref in $cUpdates.Keys) {
    $ws.Range($ref).Value = $cUpdates[$ref]
}

foreach ($ref in $dUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $dUpdates[$ref]
}

foreach ($ref in $eUpdates.Keys) {
    $ws.Range($ref).Value = $eUpdates[$ref]
}
